# Append the 2025-10-06 "8時" row (row 68) to the daily tracking sheet.
#
# A68/B68 repeat the same literal date/weekday text already used for the
# other 2025/10/06 rows (A66:A67 / B66:B67), so we copy one of those cells
# instead of assigning a literal "2025/10/06" string to .Value - the latter
# gets auto-parsed into an Excel date serial by the COM value setter, which
# is not what the source data (plain text "yyyy/mm/dd") represents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Copy($ws.Range("A68"))
$ws.Range("B67").Copy($ws.Range("B68"))
$ws.Range("C68").Value = 8
$ws.Range("D68").Value = 6
